# 7.8 History Card & Advanced Story
# Updates a handful of dialogue lines in the Kong/Lee interrogation sheet:
#  - the "missing sabers" clue line is reworded and wrapped as a green,
#    parenthesised "thought" line (rich-text style markup used by the game)
#  - three lines get their single em dash turned into a double em dash
#  - "Steward He" is renamed to "Butler He"
# Also restores the UI selection to B20 (matching the saved view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Order matters: new/changed shared strings are appended to the shared
# string table in the order they are (re)written, so we apply the edits
# in the same order the author's saved file shows them landing at the
# end of xl/sharedStrings.xml.
$ws.Range("B22").Value = " <color=#00CC00>(It’s unlikely that stealing weapons was the killer’s motive, but the missing sabers are certainly an important clue.)</color>"
$ws.Range("B3").Value = "Huh? Don’t suspect me——I injured my leg! There’s no way I could be the killer!"
$ws.Range("B5").Value = "You arrived at the banquet hall early and didn’t leave at all during the dinner——I can vouch for that."
$ws.Range("B6").Value = "That’s a relief. Ask anything you like——I’ll answer truthfully."
$ws.Range("B12").Value = "After Butler He came to inform me of the time and place, I headed out immediately."

# Restore the saved selection/view state.
$ws.Range("B20").Select() | Out-Null
